$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update row 2 (H row) values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 188
$wsOff.Range("C2").Value = 133
$wsOff.Range("D2").Value = 45
$wsOff.Range("E2").Value = 15

# Sheet "DEF" - update row 2 (H row) values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 164
$wsDef.Range("C2").Value = 102
